# Auto-generated Excel COM-interop script to update cryptos list data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "60.270.88"
$ws.Range("E2").Value = "  +4.06%  "
$ws.Range("D3").Value = "2.453.78"
$ws.Range("E3").Value = "  +4.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.99"
$ws.Range("E5").Value = "  +3.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.07"
$ws.Range("E6").Value = "  +2.42%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  +1.51%  "
$ws.Range("E9").Value = "  +5.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.83"
$ws.Range("E10").Value = "  +5.20%  "
$ws.Range("E11").Value = "  +2.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.01"
$ws.Range("E13").Value = "  +5.20%  "
$ws.Range("D14").Value = "2.884.68"
$ws.Range("E14").Value = "  +4.30%  "
$ws.Range("D15").Value = "60.141.09"
$ws.Range("E15").Value = "  +3.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000141"
$ws.Range("E16").Value = "  +5.44%  "
$ws.Range("D17").Value = "2.447.42"
$ws.Range("E17").Value = "  +3.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.51"
$ws.Range("E18").Value = "  +7.74%  "
$ws.Range("E19").Value = "  +4.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "336.06"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("E21").Value = "  +2.27%  "
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.80"
$ws.Range("E23").Value = "  +3.11%  "
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("E27").Value = "  +0.89%  "
$ws.Range("D28").Value = "0.0₃0799"
$ws.Range("E28").Value = "  +8.61%  "
$ws.Range("E29").Value = "  +4.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.34"
$ws.Range("E30").Value = "  +3.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.29"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.87"
$ws.Range("E32").Value = "  +2.30%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  +6.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.30"
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.64"
$ws.Range("E38").Value = "  +0.60%  "
$ws.Range("E39").Value = "  +2.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.419"
$ws.Range("E40").Value = "  +10.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "317.84"
$ws.Range("E41").Value = "  +8.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.76"
$ws.Range("E42").Value = "  +2.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "144.20"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("E44").Value = "  +1.99%  "
$ws.Range("E45").Value = "  +4.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.53"
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.577"
$ws.Range("E47").Value = "  +2.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.405"
$ws.Range("E48").Value = "  +5.12%  "
$ws.Range("E49").Value = "  +3.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.05"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("E51").Value = "  +5.25%  "
